# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp label update
$ws.Range("A1").Value = "Datos actualizados a 19 de Mayo de 2020 a las 16:05"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 1552304
$ws.Range("C4").Value = 2010
$ws.Range("D4").Value = 358918
$ws.Range("E4").Value = 1101314
$ws.Range("G4").Value = 91
$ws.Range("H4").Value = 92072

# --- Row 11: Alemania ---
$ws.Range("B11").Value = 177387
$ws.Range("C11").Value = 98
$ws.Range("E11").Value = 13556
$ws.Range("G11").Value = 8
$ws.Range("H11").Value = 8131

# --- Row 14: India ---
$ws.Range("B14").Value = 102335
$ws.Range("C14").Value = 2007
$ws.Range("D14").Value = 39674
$ws.Range("E14").Value = 59492

# --- Row 28: Suiza ---
$ws.Range("E28").Value = 1128
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 1890

# --- Row 60: Kazajistan ---
$ws.Range("D60").Value = 3572
$ws.Range("E60").Value = 3144

# --- Row 63: Moldavia ---
$ws.Range("D63").Value = 2508
$ws.Range("E63").Value = 3411
$ws.Range("G63").Value = 2
$ws.Range("H63").Value = 219

# --- Row 65: Oman ---
$ws.Range("E65").Value = 4070
$ws.Range("G65").Value = 2
$ws.Range("H65").Value = 27

# --- Rows 136-142: data refreshed and re-sorted by "Casos totales" (column B) descending.
# "Republica de Africa Central" now sorts in above "Etiopia", pushing the block down by one
# position (old row 142 "Republica de Africa Central" disappears from the bottom of the
# block and the other six countries each shift down one row).
$ws.Range("A136").Value = "Republica de Africa Central"
$ws.Range("B136").Value = 366
$ws.Range("C136").Value = 39
$ws.Range("D136").Value = 18
$ws.Range("E136").Value = 348
$ws.Range("G136").Value = 0
$ws.Range("H136").Value = 0

$ws.Range("A137").Value = "Etiopia"
$ws.Range("B137").Value = 365
$ws.Range("C137").Value = 13
$ws.Range("D137").Value = 120
$ws.Range("E137").Value = 240
$ws.Range("G137").Value = 0
$ws.Range("H137").Value = 5

$ws.Range("A138").Value = "Benin"
$ws.Range("B138").Value = 339
$ws.Range("C138").Value = 0
$ws.Range("D138").Value = 83
$ws.Range("E138").Value = 254
$ws.Range("G138").Value = 0
$ws.Range("H138").Value = 2

$ws.Range("A139").Value = "Isla de Man"
$ws.Range("B139").Value = 335
$ws.Range("C139").Value = 0
$ws.Range("D139").Value = 296
$ws.Range("E139").Value = 15
$ws.Range("G139").Value = 0
$ws.Range("H139").Value = 24

$ws.Range("A140").Value = "Mauricio"
$ws.Range("B140").Value = 332
$ws.Range("C140").Value = 0
$ws.Range("D140").Value = 322
$ws.Range("E140").Value = 0
$ws.Range("G140").Value = 0
$ws.Range("H140").Value = 10

$ws.Range("A141").Value = "Togo"
$ws.Range("B141").Value = 330
$ws.Range("C141").Value = 0
$ws.Range("D141").Value = 106
$ws.Range("E141").Value = 212
$ws.Range("G141").Value = 0
$ws.Range("H141").Value = 12

$ws.Range("A142").Value = "Cabo Verde"
$ws.Range("B142").Value = 328
$ws.Range("C142").Value = 0
$ws.Range("D142").Value = 85
$ws.Range("E142").Value = 240
$ws.Range("G142").Value = 0
$ws.Range("H142").Value = 3
